$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.553100000000006
$ws.Range("C3").Value = -12.88199999999999
$ws.Range("E3").Value = 16.5748
$ws.Range("E6").Value = 16.3754
$ws.Range("D8").Value = -8.876099999999997
$ws.Range("E10").Value = 16.49130000000001
$ws.Range("D11").Value = -7.394000000000003
$ws.Range("A12").Value = -21.61680000000001
$ws.Range("B14").Value = 6.2653
$ws.Range("D14").Value = -7.543400000000004
$ws.Range("D15").Value = -8.310599999999996
$ws.Range("D17").Value = -8.450699999999992
$ws.Range("C20").Value = -12.80400000000001
$ws.Range("C25").Value = -12.2782
$ws.Range("B26").Value = 4.098800000000002
$ws.Range("D26").Value = -8.665100000000008
$ws.Range("A27").Value = -22.10369999999999
$ws.Range("E27").Value = 16.6028
$ws.Range("C30").Value = -13.979
$ws.Range("B31").Value = 4.060699999999998
$ws.Range("A32").Value = -21.61020000000002
$ws.Range("E33").Value = 17.17570000000002
$ws.Range("B35").Value = 8.856700000000005
$ws.Range("A36").Value = -19.3963
$ws.Range("D36").Value = -7.494900000000005
$ws.Range("B37").Value = 8.7986
$ws.Range("A38").Value = -19.1563
$ws.Range("E39").Value = 15.7164
$ws.Range("C44").Value = -13.83249999999999
$ws.Range("B45").Value = 5.892799999999998
$ws.Range("A46").Value = -21.6444
$ws.Range("C47").Value = -12.28419999999999
$ws.Range("E47").Value = 16.5377
$ws.Range("B52").Value = 5.321500000000001
$ws.Range("A54").Value = -22.091
$ws.Range("E54").Value = 16.8417
$ws.Range("A55").Value = -22.10569999999999
$ws.Range("A56").Value = -21.58359999999999
$ws.Range("E56").Value = 16.28380000000001
$ws.Range("B57").Value = 4.902199999999994
$ws.Range("C58").Value = -13.3239
$ws.Range("E58").Value = 16.21990000000001
$ws.Range("D64").Value = -7.4154
$ws.Range("E66").Value = 17.34510000000002
$ws.Range("A67").Value = -21.55659999999998
$ws.Range("A69").Value = -21.65759999999998
$ws.Range("E69").Value = 17.38470000000002
$ws.Range("A72").Value = -21.8191
$ws.Range("E72").Value = 17.07819999999998
$ws.Range("C78").Value = -11.3475
$ws.Range("D79").Value = -6.1984
$ws.Range("E80").Value = 16.71190000000001
$ws.Range("B81").Value = 5.922399999999999
$ws.Range("E82").Value = 16.7279
$ws.Range("A83").Value = -22.0775
$ws.Range("B83").Value = 5.535999999999999
$ws.Range("E83").Value = 16.6347
$ws.Range("C84").Value = -13.95579999999999
$ws.Range("A86").Value = -22.069
$ws.Range("C89").Value = -11.171
$ws.Range("D89").Value = -5.840400000000002
$ws.Range("A91").Value = -21.48310000000001
$ws.Range("C91").Value = -11.1273
$ws.Range("C92").Value = -11.4153
$ws.Range("A93").Value = -21.2743
$ws.Range("C96").Value = -13.4816
$ws.Range("A99").Value = -20.15639999999998
$ws.Range("B100").Value = 4.760800000000001
$ws.Range("B102").Value = 8.136100000000001
$ws.Range("C102").Value = -14.08880000000001
